$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '25.907.89'
$ws.Range('E2').Value = '  -0.14%  '
$ws.Range('D3').Value = '1.585.17'
$ws.Range('E3').Value = '  -2.01%  '
$ws.Range('E4').Value = '  -0.35%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '210.36'
$ws.Range('E5').Value = '  -0.64%  '
$ws.Range('E6').Value = '  -0.27%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.479'
$ws.Range('E7').Value = '  -1.60%  '
$ws.Range('E8').Value = '  +0.09%  '
$ws.Range('E9').Value = '  -1.14%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '18.10'
$ws.Range('E10').Value = '  -0.47%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.0791'
$ws.Range('E11').Value = '  +0.09%  '
$ws.Range('D12').Value = '1.804.22'
$ws.Range('E12').Value = '  -2.17%  '
$ws.Range('D13').Value = '1.589.82'
$ws.Range('E13').Value = '  -1.31%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '4.03'
$ws.Range('E14').Value = '  -2.38%  '
$ws.Range('E15').Value = '  -2.04%  '
$ws.Range('D16').Value = '25.878.29'
$ws.Range('E16').Value = '  -0.31%  '
$ws.Range('D17').Value = '0.0₃0725'
$ws.Range('E17').Value = '  -1.10%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '59.85'
$ws.Range('E18').Value = '  -2.82%  '
$ws.Range('E19').Value = '  -0.25%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '192.53'
$ws.Range('E20').Value = '  +0.73%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '4.19'
$ws.Range('E21').Value = '  -0.98%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '9.36'
$ws.Range('E22').Value = '  -1.23%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '5.95'
$ws.Range('E23').Value = '  -0.95%  '
$ws.Range('E24').Value = '  +0.36%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '141.36'
$ws.Range('E25').Value = '  -1.81%  '
$ws.Range('E26').Value = '  -0.38%  '
$ws.Range('E27').Value = '  -0.72%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '15.09'
$ws.Range('E28').Value = '  -0.62%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '6.45'
$ws.Range('E29').Value = '  -2.51%  '
$ws.Range('E30').Value = '  -5.30%  '
$ws.Range('E31').Value = '  -0.92%  '
$ws.Range('E32').Value = '  +0.37%  '
$ws.Range('E33').Value = '  -1.77%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '1.50'
$ws.Range('E34').Value = '  +0.67%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '2.36'
$ws.Range('E35').Value = '  -1.96%  '
$ws.Range('D36').Value = '1.096.25'
$ws.Range('E36').Value = '  -2.68%  '
$ws.Range('B37').Value = 'VeChain'
$ws.Range('C37').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '0.0154'
$ws.Range('E37').Value = '  +0.69%  '
$ws.Range('B38').Value = 'PaxDollar'
$ws.Range('C38').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '1.00'
$ws.Range('E38').Value = '  -0.38%  '
$ws.Range('B39').Value = 'MXToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '2.35'
$ws.Range('E39').Value = '  -2.49%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.501'
$ws.Range('E40').Value = '  -2.82%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.777'
$ws.Range('E41').Value = '  -5.54%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '0.806'
$ws.Range('E42').Value = '  +6.83%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '93.48'
$ws.Range('E43').Value = '  -4.21%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '5.14'
$ws.Range('E44').Value = '  +0.60%  '
$ws.Range('D45').Value = '1.717.95'
$ws.Range('E45').Value = '  -2.15%  '
$ws.Range('D46').Value = '0.0₆0113'
$ws.Range('E46').Value = '  -0.61%  '
$ws.Range('E47').Value = '  +1.32%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '53.13'
$ws.Range('E48').Value = '  -1.31%  '
$ws.Range('E49').Value = '  -0.95%  '
$ws.Range('E50').Value = '  -0.93%  '
$ws.Range('E51').Value = '  -0.32%  '
